# Apply edits to testdata.xlsx
$wb = $excel.ActiveWorkbook

# --- InvalidLoginData sheet (2nd sheet) ---
$ws2 = $wb.Worksheets.Item(2)

# Clear the old B4 value and the old row 7 data
$ws2.Range("B4").ClearContents()
$ws2.Range("A7").ClearContents()
$ws2.Range("B7").ClearContents()

# Set the new values for rows 4 and 5
$ws2.Range("A4").Value = "testuser@example.com"
$ws2.Range("A5").Value = "test+user@example.com"
$ws2.Range("B5").Value = "Test@123"

# Update selection on InvalidLoginData sheet and make it the active sheet/tab
$ws2.Select()
$ws2.Range("A4:XFD4").Select()
